# Auto-generated Excel COM-interop script
# Applies numeric corrections to the FFXIV leve-profit workbook ("Phantom_Profits")
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR as per the target diff.

$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3195.6667
$ws.Range("I40").Value = 2793.375
$ws.Range("J40").Value = 4000.25
$ws.Range("K40").Value = 2793.375
$ws.Range("L40").Value = 4000.25
$ws.Range("M40").Value = -2618.375
$ws.Range("N40").Value = -4350.25
$ws.Range("H80").Value = 4065.5
$ws.Range("I80").Value = 4096.2856
$ws.Range("J80").Value = 4022.4
$ws.Range("K80").Value = 12288.8568
$ws.Range("L80").Value = 12067.2
$ws.Range("M80").Value = -11290.8568
$ws.Range("N80").Value = -14063.2
$ws.Range("H81").Value = 43000
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H83").Value = 4065.5
$ws.Range("I83").Value = 4096.2856
$ws.Range("J83").Value = 4022.4
$ws.Range("K83").Value = 36866.5704
$ws.Range("L83").Value = 36201.6
$ws.Range("M83").Value = -31874.5704
$ws.Range("N83").Value = -46185.6
$ws.Range("H84").Value = 43000
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H130").Value = 27744
$ws.Range("J130").Value = 25779
$ws.Range("L130").Value = 25779
$ws.Range("N130").Value = -35819
$ws.Range("H132").Value = 4680.2607
$ws.Range("I132").Value = 4632.75
$ws.Range("K132").Value = 13898.25
$ws.Range("M132").Value = -11368.25
$ws.Range("H137").Value = 1976.8182
$ws.Range("J137").Value = 2020.75
$ws.Range("L137").Value = 6062.25
$ws.Range("N137").Value = -11162.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1201.7858
$ws.Range("I32").Value = 1201.7858
$ws.Range("K32").Value = 1201.7858
$ws.Range("M32").Value = -914.7858000000001
$ws.Range("H45").Value = 2399.5
$ws.Range("I45").Value = 2000
$ws.Range("J45").Value = 2799
$ws.Range("K45").Value = 2000
$ws.Range("L45").Value = 2799
$ws.Range("M45").Value = -1623
$ws.Range("N45").Value = -3553
$ws.Range("H95").Value = 51529.145
$ws.Range("J95").Value = 51529.145
$ws.Range("L95").Value = 51529.145
$ws.Range("N95").Value = -57021.145
$ws.Range("H110").Value = 5111.636
$ws.Range("I110").Value = 5322.8
$ws.Range("K110").Value = 5322.8
$ws.Range("M110").Value = -3277.8
$ws.Range("H122").Value = 1642.6
$ws.Range("I122").Value = 1617.0714
$ws.Range("K122").Value = 4851.2142
$ws.Range("M122").Value = -2401.2142
$ws.Range("H132").Value = 3590.5454
$ws.Range("I132").Value = 5399.8
$ws.Range("J132").Value = 2082.8333
$ws.Range("K132").Value = 16199.4
$ws.Range("L132").Value = 6248.499899999999
$ws.Range("M132").Value = -13669.4
$ws.Range("N132").Value = -11308.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H134").Value = 1328.8
$ws.Range("I134").Value = 1328.8
$ws.Range("K134").Value = 3986.4
$ws.Range("M134").Value = -1451.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H100").Value = 49999
$ws.Range("J100").Value = 49999
$ws.Range("L100").Value = 49999
$ws.Range("N100").Value = -52163
$ws.Range("H107").Value = 787.875
$ws.Range("I107").Value = 814.8570999999999
$ws.Range("K107").Value = 814.8570999999999
$ws.Range("M107").Value = 1105.1429
$ws.Range("H112").Value = 40102
$ws.Range("J112").Value = 40102
$ws.Range("L112").Value = 40102
$ws.Range("N112").Value = -43056
$ws.Range("H132").Value = 4999
$ws.Range("J132").Value = 4999
$ws.Range("L132").Value = 14997
$ws.Range("N132").Value = -20057
$ws.Range("H134").Value = 2724.5
$ws.Range("I134").Value = 2716.95
$ws.Range("K134").Value = 8150.849999999999
$ws.Range("M134").Value = -5615.849999999999
$ws.Range("H138").Value = 179999
$ws.Range("J138").Value = 179999
$ws.Range("L138").Value = 179999
$ws.Range("N138").Value = -190279

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5000
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H23").Value = 218.6
$ws.Range("J23").Value = 246
$ws.Range("L23").Value = 738
$ws.Range("N23").Value = -1208
$ws.Range("H62").Value = 5665
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 5665
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H92").Value = 286.75
$ws.Range("I92").Value = 286.75
$ws.Range("K92").Value = 860.25
$ws.Range("M92").Value = 387.75
$ws.Range("H97").Value = 366.8
$ws.Range("I97").Value = 403
$ws.Range("J97").Value = 222
$ws.Range("K97").Value = 1209
$ws.Range("L97").Value = 666
$ws.Range("M97").Value = -713
$ws.Range("N97").Value = -1658
$ws.Range("H98").Value = 386
$ws.Range("J98").Value = 547.5
$ws.Range("L98").Value = 1642.5
$ws.Range("N98").Value = -4638.5
$ws.Range("H135").Value = 5000
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H138").Value = 13553.4
$ws.Range("I138").Value = 11961.143
$ws.Range("K138").Value = 35883.429
$ws.Range("M138").Value = -30743.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 334333.34
$ws.Range("I7").Value = 500500
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 500500
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = -500388
$ws.Range("N7").Value = -2224
$ws.Range("H8").Value = 334333.34
$ws.Range("I8").Value = 500500
$ws.Range("J8").Value = 2000
$ws.Range("K8").Value = 500500
$ws.Range("L8").Value = 2000
$ws.Range("M8").Value = -500361
$ws.Range("N8").Value = -2278
$ws.Range("H104").Value = 49912.5
$ws.Range("J104").Value = 49912.5
$ws.Range("L104").Value = 49912.5
$ws.Range("N104").Value = -56900.5
$ws.Range("H116").Value = 150000
$ws.Range("J116").Value = 150000
$ws.Range("L116").Value = 150000
$ws.Range("N116").Value = -159178
$ws.Range("H132").Value = 2903.7144
$ws.Range("I132").Value = 2953.682
$ws.Range("K132").Value = 8861.045999999998
$ws.Range("M132").Value = -6331.045999999998
$ws.Range("H135").Value = 77333.336
$ws.Range("J135").Value = 77333.336
$ws.Range("L135").Value = 77333.336
$ws.Range("N135").Value = -87473.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2072.6667
$ws.Range("I7").Value = 1750.6666
$ws.Range("J7").Value = 2716.6667
$ws.Range("K7").Value = 1750.6666
$ws.Range("L7").Value = 2716.6667
$ws.Range("M7").Value = -1638.6666
$ws.Range("N7").Value = -2940.6667
$ws.Range("H9").Value = 315.2
$ws.Range("J9").Value = 199
$ws.Range("L9").Value = 199
$ws.Range("N9").Value = -647
$ws.Range("H22").Value = 2750
$ws.Range("I22").Value = 1500
$ws.Range("J22").Value = 3375
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 3375
$ws.Range("M22").Value = -1205
$ws.Range("N22").Value = -3965
$ws.Range("H27").Value = 2750
$ws.Range("I27").Value = 1500
$ws.Range("J27").Value = 3375
$ws.Range("K27").Value = 1500
$ws.Range("L27").Value = 3375
$ws.Range("M27").Value = -1393
$ws.Range("N27").Value = -3589
$ws.Range("H40").Value = 10887.667
$ws.Range("I40").Value = 5998.625
$ws.Range("K40").Value = 5998.625
$ws.Range("M40").Value = -5862.625
$ws.Range("H68").Value = 7257.1665
$ws.Range("I68").Value = 6708.6
$ws.Range("K68").Value = 6708.6
$ws.Range("M68").Value = -5959.6
$ws.Range("H71").Value = 7257.1665
$ws.Range("I71").Value = 6708.6
$ws.Range("K71").Value = 33543
$ws.Range("M71").Value = -29799
$ws.Range("H126").Value = 2072.6667
$ws.Range("I126").Value = 1750.6666
$ws.Range("J126").Value = 2716.6667
$ws.Range("K126").Value = 5251.9998
$ws.Range("L126").Value = 8150.000100000001
$ws.Range("M126").Value = -2781.9998
$ws.Range("N126").Value = -13090.0001
$ws.Range("H132").Value = 2799.8333
$ws.Range("I132").Value = 2334.6667
$ws.Range("K132").Value = 7004.000100000001
$ws.Range("M132").Value = -4474.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 3000
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H132").Value = 1339.5714
$ws.Range("I132").Value = 1339.5714
$ws.Range("K132").Value = 4018.7142
$ws.Range("M132").Value = -1488.7142
